$d = $word.ActiveDocument
$d.Content.Find.Execute("09/05/2014", $true, $false, $false, $false, $false,
                         $true, 1, $false, "03/05/2014", 2)
